$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Simple numeric updates (no type/style change) ---
$ws.Range("L15").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 39
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -15.217391304347
$ws.Range("L16").Value = -4.878048780487
$ws.Range("M16").Value = 14.705882352941
$ws.Range("N16").Value = -82.894736842105
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 35.294117647058
$ws.Range("L17").Value = 119.047619047619
$ws.Range("M17").Value = 64.285714285714
$ws.Range("N17").Value = 15
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 41.666666666666
$ws.Range("M18").Value = 18.604651162790
$ws.Range("N18").Value = -83.168316831683
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 88.888888888888
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = 11.290322580645
$ws.Range("I19").Value = 278
$ws.Range("J19").Value = 280
$ws.Range("K19").Value = -0.714285714285
$ws.Range("L19").Value = 22.466960352422
$ws.Range("M19").Value = 7.335907335907
$ws.Range("N19").Value = -63.707571801566
$ws.Range("C20").Value = 5
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = 61.904761904761
$ws.Range("L20").Value = 41.666666666666
$ws.Range("M20").Value = 209.090909090909
$ws.Range("N20").Value = -92.543859649122
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 60
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 10.752688172043
$ws.Range("I21").Value = 454
$ws.Range("J21").Value = 436
$ws.Range("K21").Value = 4.128440366972
$ws.Range("L21").Value = 28.248587570621
$ws.Range("M21").Value = 19.473684210526
$ws.Range("N21").Value = -74.819744869661
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 19
$ws.Range("K23").Value = 280
$ws.Range("L23").Value = 90
$ws.Range("M23").Value = 111.111111111111
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -61.111111111111
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -30.357142857142
$ws.Range("I24").Value = 397
$ws.Range("J24").Value = 543
$ws.Range("K24").Value = -26.887661141804
$ws.Range("L24").Value = -29.233511586452
$ws.Range("M24").Value = 6.434316353887
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 20
$ws.Range("I25").Value = 89
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = 5.952380952380
$ws.Range("L25").Value = 43.548387096774
$ws.Range("M25").Value = -9.183673469387
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = 10
$ws.Range("K26").Value = 100
$ws.Range("L26").Value = 42.857142857142
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -26.315789473684
$ws.Range("L27").Value = 16.666666666666

# --- Numeric -> Text (shared string) conversions ---
$ws.Range("C17").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "'***.*"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "'***.*"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "'***.*"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null

# --- Text (shared string) -> Numeric conversions ---
$ws.Range("D18").Value = 2
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = -100
$ws.Range("M14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
